$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Tables": a service-type column header is added to the second
# "pff_api_servicelog" table block (G3 = "servicetype(text)" - reuses the
# existing shared string already used at K4).
# ---------------------------------------------------------------------------
$wsTables = $wb.Worksheets.Item("Tables")
$wsTables.Range("G3").Value = "servicetype(text)"

# ---------------------------------------------------------------------------
# Sheet "Jar files": a second jar entry is appended.
# ---------------------------------------------------------------------------
$wsJar = $wb.Worksheets.Item("Jar files")
$wsJar.Range("A3").Value = 2
$wsJar.Range("B3").Value = "dockerized upload file jar for token api"
$wsJar.Range("C3").Value = 8081
$wsJar.Columns.Item(1).AutoFit() | Out-Null
$wsJar.Columns.Item(2).AutoFit() | Out-Null
$wsJar.Columns.Item(3).AutoFit() | Out-Null
$wsJar.Range("B3").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "Plugins": a fourth plugin entry ("prometheus") is appended.
# ---------------------------------------------------------------------------
$wsPlugins = $wb.Worksheets.Item("Plugins")
$wsPlugins.Range("A5").Value = 4
$wsPlugins.Range("B5").Value = "prometheus"
$wsPlugins.Range("B7").Select() | Out-Null

# ---------------------------------------------------------------------------
# Leave the workbook the way it was originally focused: "Tables" sheet
# active, cell C4 selected.
# ---------------------------------------------------------------------------
$wsTables.Select() | Out-Null
$wsTables.Range("C4").Select() | Out-Null
